# Vtn-Itgb6.xlsx: refresh with the newly-computed TPM-based NATMI output.
# The "Inflammatory-Mac" target cluster no longer appears in this run, so
# its three rows are removed (also drops the now-unused shared string),
# and every remaining metric cell is rewritten with the recalculated value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so the row numbers of rows above stay valid.
$ws.Rows(14).Delete()   # MuSCs  -> Inflammatory-Mac
$ws.Rows(9).Delete()    # FAPs   -> Inflammatory-Mac
$ws.Rows(4).Delete()    # ECs    -> Inflammatory-Mac

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.440985666666667
$ws.Cells.Item(2, 8).Value = 4.322957000000001
$ws.Cells.Item(2, 9).Value = 0.1098365531732288
$ws.Cells.Item(2, 10).Value = 0.1230162332390494
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1861273333333333
$ws.Cells.Item(2, 14).Value = 0.5583819999999999
$ws.Cells.Item(2, 15).Value = 0.01396039250968796
$ws.Cells.Item(2, 16).Value = 0.01786028945924599
$ws.Cells.Item(2, 17).Value = 0.2682068195082222
$ws.Cells.Item(2, 18).Value = 2.413861375574
$ws.Cells.Item(2, 19).Value = 0.001533361394209486
$ws.Cells.Item(2, 20).Value = 0.00219710553383554

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.440985666666667
$ws.Cells.Item(3, 8).Value = 4.322957000000001
$ws.Cells.Item(3, 9).Value = 0.1098365531732288
$ws.Cells.Item(3, 10).Value = 0.1230162332390494
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.379107333333333
$ws.Cells.Item(3, 14).Value = 13.137322
$ws.Cells.Item(3, 15).Value = 0.3284528721308331
$ws.Cells.Item(3, 16).Value = 0.4202076242416848
$ws.Cells.Item(3, 17).Value = 6.310230900128222
$ws.Cells.Item(3, 18).Value = 56.792078101154
$ws.Cells.Item(3, 19).Value = 0.03607613135469796
$ws.Cells.Item(3, 20).Value = 0.05169235911254191

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.440985666666667
$ws.Cells.Item(4, 8).Value = 4.322957000000001
$ws.Cells.Item(4, 9).Value = 0.1098365531732288
$ws.Cells.Item(4, 10).Value = 0.1230162332390494
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 8.733703
$ws.Cells.Item(4, 14).Value = 17.467406
$ws.Cells.Item(4, 15).Value = 0.655067258309039
$ws.Cells.Item(4, 16).Value = 0.5587087822712233
$ws.Cells.Item(4, 17).Value = 12.58514083992367
$ws.Cells.Item(4, 18).Value = 75.51084503954202
$ws.Cells.Item(4, 19).Value = 0.07195032974930196
$ws.Cells.Item(4, 20).Value = 0.06873024987258206

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 1.440985666666667
$ws.Cells.Item(5, 8).Value = 4.322957000000001
$ws.Cells.Item(5, 9).Value = 0.1098365531732288
$ws.Cells.Item(5, 10).Value = 0.1230162332390494
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.033591
$ws.Cells.Item(5, 14).Value = 0.100773
$ws.Cells.Item(5, 15).Value = 0.002519477050439994
$ws.Cells.Item(5, 16).Value = 0.003223304027845805
$ws.Cells.Item(5, 17).Value = 0.04840414952900001
$ws.Cells.Item(5, 18).Value = 0.435637345761
$ws.Cells.Item(5, 19).Value = 0.0002767306750193821
$ws.Cells.Item(5, 20).Value = 0.0003965187200898469

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.461641333333333
$ws.Cells.Item(6, 8).Value = 22.384924
$ws.Cells.Item(6, 9).Value = 0.5687502547919595
$ws.Cells.Item(6, 10).Value = 0.6369966279614609
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1861273333333333
$ws.Cells.Item(6, 14).Value = 0.5583819999999999
$ws.Cells.Item(6, 15).Value = 0.01396039250968796
$ws.Cells.Item(6, 16).Value = 0.01786028945924599
$ws.Cells.Item(6, 17).Value = 1.388815403663111
$ws.Cells.Item(6, 18).Value = 12.499338632968
$ws.Cells.Item(6, 19).Value = 0.00793997679688079
$ws.Cells.Item(6, 20).Value = 0.01137694415995532

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.461641333333333
$ws.Cells.Item(7, 8).Value = 22.384924
$ws.Cells.Item(7, 9).Value = 0.5687502547919595
$ws.Cells.Item(7, 10).Value = 0.6369966279614609
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.379107333333333
$ws.Cells.Item(7, 14).Value = 13.137322
$ws.Cells.Item(7, 15).Value = 0.3284528721308331
$ws.Cells.Item(7, 16).Value = 0.4202076242416848
$ws.Cells.Item(7, 17).Value = 32.67532828150311
$ws.Cells.Item(7, 18).Value = 294.077954533528
$ws.Cells.Item(7, 19).Value = 0.1868076547115622
$ws.Cells.Item(7, 20).Value = 0.2676708396856499

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.461641333333333
$ws.Cells.Item(8, 8).Value = 22.384924
$ws.Cells.Item(8, 9).Value = 0.5687502547919595
$ws.Cells.Item(8, 10).Value = 0.6369966279614609
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.733703
$ws.Cells.Item(8, 14).Value = 17.467406
$ws.Cells.Item(8, 15).Value = 0.655067258309039
$ws.Cells.Item(8, 16).Value = 0.5587087822712233
$ws.Cells.Item(8, 17).Value = 65.16775929785733
$ws.Cells.Item(8, 18).Value = 391.006555787144
$ws.Cells.Item(8, 19).Value = 0.3725696700691363
$ws.Cells.Item(8, 20).Value = 0.3558956103192233

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.461641333333333
$ws.Cells.Item(9, 8).Value = 22.384924
$ws.Cells.Item(9, 9).Value = 0.5687502547919595
$ws.Cells.Item(9, 10).Value = 0.6369966279614609
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.033591
$ws.Cells.Item(9, 14).Value = 0.100773
$ws.Cells.Item(9, 15).Value = 0.002519477050439994
$ws.Cells.Item(9, 16).Value = 0.003223304027845805
$ws.Cells.Item(9, 17).Value = 0.250643994028
$ws.Cells.Item(9, 18).Value = 2.255795946252
$ws.Cells.Item(9, 19).Value = 0.001432953214380241
$ws.Cells.Item(9, 20).Value = 0.002053233796632373

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.2167365
$ws.Cells.Item(10, 8).Value = 8.433472999999999
$ws.Cells.Item(10, 9).Value = 0.3214131920348118
$ws.Cells.Item(10, 10).Value = 0.2399871387994896
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1861273333333333
$ws.Cells.Item(10, 14).Value = 0.5583819999999999
$ws.Cells.Item(10, 15).Value = 0.01396039250968796
$ws.Cells.Item(10, 16).Value = 0.01786028945924599
$ws.Cells.Item(10, 17).Value = 0.7848499201143332
$ws.Cells.Item(10, 18).Value = 4.709099520685999
$ws.Cells.Item(10, 19).Value = 0.004487054318597684
$ws.Cells.Item(10, 20).Value = 0.004286239765455129

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.2167365
$ws.Cells.Item(11, 8).Value = 8.433472999999999
$ws.Cells.Item(11, 9).Value = 0.3214131920348118
$ws.Cells.Item(11, 10).Value = 0.2399871387994896
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 4.379107333333333
$ws.Cells.Item(11, 14).Value = 13.137322
$ws.Cells.Item(11, 15).Value = 0.3284528721308331
$ws.Cells.Item(11, 16).Value = 0.4202076242416848
$ws.Cells.Item(11, 17).Value = 18.46554172988433
$ws.Cells.Item(11, 18).Value = 110.793250379306
$ws.Cells.Item(11, 19).Value = 0.1055690860645729
$ws.Cells.Item(11, 20).Value = 0.100844425443493

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.2167365
$ws.Cells.Item(12, 8).Value = 8.433472999999999
$ws.Cells.Item(12, 9).Value = 0.3214131920348118
$ws.Cells.Item(12, 10).Value = 0.2399871387994896
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 8.733703
$ws.Cells.Item(12, 14).Value = 17.467406
$ws.Cells.Item(12, 15).Value = 0.655067258309039
$ws.Cells.Item(12, 16).Value = 0.5587087822712233
$ws.Cells.Item(12, 17).Value = 36.8277242202595
$ws.Cells.Item(12, 18).Value = 147.310896881038
$ws.Cells.Item(12, 19).Value = 0.2105472584906008
$ws.Cells.Item(12, 20).Value = 0.1340829220794179

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.2167365
$ws.Cells.Item(13, 8).Value = 8.433472999999999
$ws.Cells.Item(13, 9).Value = 0.3214131920348118
$ws.Cells.Item(13, 10).Value = 0.2399871387994896
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.033591
$ws.Cells.Item(13, 14).Value = 0.100773
$ws.Cells.Item(13, 15).Value = 0.002519477050439994
$ws.Cells.Item(13, 16).Value = 0.003223304027845805
$ws.Cells.Item(13, 17).Value = 0.1416443957715
$ws.Cells.Item(13, 18).Value = 0.8498663746289999
$ws.Cells.Item(13, 19).Value = 0.000809793161040371
$ws.Cells.Item(13, 20).Value = 0.0007735515111235851
